$wb = $excel.ActiveWorkbook

# The "BUNDLE-NAME" example file name in the bitstream-metadata sheet is
# replaced with a placeholder that the integration test substitutes at
# run time (DSC-723: BulkImport allows attaching any local file).
$ws = $wb.Worksheets.Item("bitstream-metadata")
$ws.Range("B2").Value = "THIS IS DYNAMICALLY REPLACED BY THE TEST"

# Restore the active selection on that sheet to B2.
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
